# Update the starting value of the "Pozicija" (position) counter in column B.
# B2 used to hold the literal value 1; it now holds the real starting
# position number used by the order-import. The remaining cells in the
# column (B3:B15) already contain the formula "=<previous row>+1" and will
# recalculate automatically from this new seed value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1224547468

# Reflect the author's final cell selection (cell H4) when the file was saved.
$ws.Range("H4").Select()
